$d = $word.ActiveDocument

# Update the date/title line
$d.Content.Find.Execute("2024-05-22 Wednesday", $true, $false, $false, $false, $false,
                         $true, 1, $false, "2024-05-23 Thursday", 2)

# Update the division problems in the single table, cell by cell, so that
# cells sharing the same before/after text across different positions do
# not collide with each other (e.g. "38÷3=" is both a source and a target).
$t = $d.Tables.Item(1)

$t.Cell(1,1).Range.Text  = "40÷6="
$t.Cell(1,2).Range.Text  = "84÷2="
$t.Cell(1,3).Range.Text  = "54÷9="
$t.Cell(1,4).Range.Text  = "23÷8="
$t.Cell(1,5).Range.Text  = "39÷6="

$t.Cell(5,1).Range.Text  = "67÷9="
$t.Cell(5,2).Range.Text  = "38÷3="
$t.Cell(5,3).Range.Text  = "61÷8="
$t.Cell(5,4).Range.Text  = "44÷9="
$t.Cell(5,5).Range.Text  = "70÷3="

$t.Cell(9,1).Range.Text  = "90÷5="
$t.Cell(9,2).Range.Text  = "59÷9="
$t.Cell(9,3).Range.Text  = "95÷5="
$t.Cell(9,4).Range.Text  = "31÷6="
$t.Cell(9,5).Range.Text  = "27÷5="

$t.Cell(13,1).Range.Text = "48÷6="
$t.Cell(13,2).Range.Text = "40÷5="
$t.Cell(13,3).Range.Text = "10÷7="
$t.Cell(13,4).Range.Text = "28÷8="
$t.Cell(13,5).Range.Text = "95÷2="

$t.Cell(17,1).Range.Text = "19÷9="
$t.Cell(17,2).Range.Text = "32÷5="
$t.Cell(17,3).Range.Text = "39÷2="
$t.Cell(17,4).Range.Text = "36÷8="
$t.Cell(17,5).Range.Text = "15÷3="
